$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'maa://24702 (94.48), maa://25390 (96.15), maa://36681 (87.34)'
$ws.Range('L2').Value = '*maa://24633 (56.52), *maa://30515 (69.9), *maa://34787 (73.33), maa://39402 (91.53), ***maa://20792 (11.93), ***maa://29083 (27.78)'
$ws.Range('AB2').Value = 'maa://21246 (91.41), maa://36684 (95.5), ***maa://22731 (6.67)'
$ws.Range('H3').Value = 'maa://21247 (98.55), *maa://22748 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.82), maa://20276 (86.44), *maa://22749 (72.73)'
$ws.Range('P3').Value = 'maa://21249 (94.42), maa://26254 (96.55)'
$ws.Range('X3').Value = 'maa://27396 (84.16), maa://27484 (96.55), maa://27480 (83.33)'
$ws.Range('D4').Value = 'maa://24632 (94.08), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (96.55), maa://27295 (85.71), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range('X4').Value = '**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (89.09), ***maa://36683 (28.26)'
$ws.Range('D5').Value = 'maa://21245 (84.81), maa://22744 (84.0)'
$ws.Range('L7').Value = 'maa://28624 (92.59), maa://24957 (97.73)'
$ws.Range('AF7').Value = '*maa://26191 (68.24), *maa://36671 (68.0), maa://45272 (100.0), *maa://42530 (62.5)'
$ws.Range('A8').Value = '更新日期：2025.02.14 13:17:46'
$ws.Range('P8').Value = 'maa://32931 (83.62), *maa://21916 (61.54), maa://23252 (91.18), maa://37496 (96.77), **maa://22759 (45.45)'
$ws.Range('X8').Value = 'maa://21411 (95.91)'
$ws.Range('AF8').Value = '*maa://24479 (78.16), *maa://21990 (51.85)'
$ws.Range('L9').Value = 'maa://22762 (92.31), *maa://39552 (75.0)'
$ws.Range('X9').Value = 'maa://26223 (97.84)'
$ws.Range('D10').Value = '***maa://25695 (18.82), ***maa://34206 (20.0), ***maa://39951 (15.69), ***maa://39243 (28.57), *maa://45271 (57.14)'
$ws.Range('H10').Value = 'maa://32651 (94.12)'
$ws.Range('X10').Value = 'maa://22301 (97.75), maa://22726 (100.0), maa://45828 (88.89)'
$ws.Range('AB11').Value = 'maa://29912 (97.14), maa://22516 (88.37), *maa://20794 (52.24)'
$ws.Range('D12').Value = 'maa://30766 (89.29), *maa://36678 (71.43)'
$ws.Range('H12').Value = 'maa://21867 (89.88), ***maa://45826 (20.0)'
$ws.Range('AB12').Value = 'maa://23669 (95.5), maa://36677 (93.1), maa://39872 (91.67)'
$ws.Range('AF12').Value = '*maa://28932 (77.4), *maa://20106 (63.96), *maa://22769 (64.29)'
$ws.Range('D13').Value = 'maa://24999 (92.04), maa://36673 (93.24), maa://25001 (85.71)'
$ws.Range('AF13').Value = '**maa://22737 (33.33), maa://39883 (91.3), *maa://39885 (53.33)'
$ws.Range('D15').Value = '*maa://22743 (77.62), maa://22734 (84.17), *maa://30808 (64.18), **maa://36048 (45.0), maa://45058 (91.67)'
$ws.Range('H15').Value = 'maa://24304 (87.91), maa://21478 (89.19)'
$ws.Range('AF15').Value = 'maa://21364 (81.1), *maa://36666 (78.5), *maa://22766 (68.64)'
$ws.Range('T16').Value = 'maa://22729 (94.94), *maa://28648 (69.12), maa://36674 (82.35)'
$ws.Range('H17').Value = 'maa://22430 (88.66), maa://39599 (86.0)'
$ws.Range('D18').Value = 'maa://24570 (97.31)'
$ws.Range('H18').Value = 'maa://24421 (88.98)'
$ws.Range('L18').Value = 'maa://22466 (90.0), *maa://22732 (51.14)'
$ws.Range('X18').Value = 'maa://21917 (96.91), maa://22741 (85.71)'
$ws.Range('AF18').Value = '*maa://24313 (59.64), **maa://29784 (46.43)'
$ws.Range('AB19').Value = '*maa://30709 (65.28), *maa://36668 (57.5)'
$ws.Range('L20').Value = 'maa://41331 (85.62)'
$ws.Range('X22').Value = 'maa://21282 (98.6), *maa://37649 (65.52)'
$ws.Range('L23').Value = 'maa://39756 (95.54), maa://39875 (94.37)'
$ws.Range('D25').Value = 'maa://29753 (95.11)'
$ws.Range('H25').Value = '*maa://29063 (73.75), *maa://25311 (73.53), ***maa://22725 (4.84), *maa://45047 (62.5)'
$ws.Range('X26').Value = 'maa://24389 (96.55)'
$ws.Range('T28').Value = 'maa://23263 (95.24), *maa://29765 (63.41)'
$ws.Range('X28').Value = 'maa://39929 (90.55), maa://41749 (90.36), ***maa://39723 (13.89)'
$ws.Range('AF28').Value = 'maa://36660 (92.15), *maa://36701 (65.52)'
$ws.Range('L29').Value = 'maa://28432 (93.41), *maa://28440 (79.44), maa://31400 (100.0), *maa://28650 (71.43)'
$ws.Range('AF29').Value = '*maa://24080 (68.85), maa://42865 (81.03), ***maa://34960 (8.33)'
$ws.Range('L31').Value = 'maa://35926 (93.38), maa://36258 (84.96), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.5), maa://36667 (98.73), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('L37').Value = 'maa://45718 (98.32), maa://45789 (100.0)'
$ws.Range('AF38').Value = 'maa://36697 (86.12)'
$ws.Range('T39').Value = 'maa://45788 (82.42), maa://45790 (81.82)'
$ws.Range('P40').Value = 'maa://23278 (95.54), maa://21386 (95.77), maa://36664 (89.29), maa://45550 (100.0)'
$ws.Range('H55').Value = 'maa://32532 (92.26)'
$ws.Range('H58').Value = '*maa://37964 (61.11)'
$ws.Range('H59').Value = 'maa://31270 (95.31), maa://27746 (82.3)'
